$d = $word.ActiveDocument

# Locate the "Unsubscribe" hyperlink and remove it completely (link + its
# display text), leaving the rest of the paragraph intact.
for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "Unsubscribe") {
        $r = $d.Range($h.Range.Start, $h.Range.End)
        $r.Delete()
    }
}

# The run immediately preceding the (now removed) hyperlink held four
# spaces as left-padding; trim it down to three spaces to match the
# updated layout.
$d.Content.Find.Execute("    ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "   ", 2) | Out-Null
